# Add "Vaccine 1 ..." and "Vaccine 2 ..." columns (CY..DH / cols 103..112)
# to the Monitorees sheet so the vaccine table can be populated on import.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column layout (1-based):
#  CY=103 Group Name        DD=108 Group Name
#  CZ=104 Product Name      DE=109 Product Name
#  DA=105 Administration Dt DF=110 Administration Dt  (text/date-as-text style)
#  DB=106 Dose Number       DG=111 Dose Number
#  DC=107 Notes             DH=112 Notes

$COL_V1_GROUP   = 103
$COL_V1_PRODUCT = 104
$COL_V1_DATE    = 105
$COL_V1_DOSE    = 106
$COL_V1_NOTES   = 107
$COL_V2_GROUP   = 108
$COL_V2_PRODUCT = 109
$COL_V2_DATE    = 110
$COL_V2_DOSE    = 111
$COL_V2_NOTES   = 112

function Set-Text($row, $col, $text) {
    $ws.Cells.Item($row, $col).Value = $text
}

function Set-DateText($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

function Set-Num($row, $col, $num) {
    $ws.Cells.Item($row, $col).Value = $num
}

# ---- Header row ----
Set-Text     1 $COL_V1_GROUP   "Vaccine 1 Group Name"
Set-Text     1 $COL_V1_PRODUCT "Vaccine 1 Product Name"
Set-DateText 1 $COL_V1_DATE    "Vaccine 1 Administration Date"
Set-Text     1 $COL_V1_DOSE    "Vaccine 1 Dose Number"
Set-Text     1 $COL_V1_NOTES   "Vaccine 1 Notes"
Set-Text     1 $COL_V2_GROUP   "Vaccine 2 Group Name"
Set-Text     1 $COL_V2_PRODUCT "Vaccine 2 Product Name"
Set-DateText 1 $COL_V2_DATE    "Vaccine 2 Administration Date"
Set-Text     1 $COL_V2_DOSE    "Vaccine 2 Dose Number"
Set-Text     1 $COL_V2_NOTES   "Vaccine 2 Notes"

# ---- Row 2 ----
Set-Text     2 $COL_V1_GROUP   "COVID-19"
Set-Text     2 $COL_V1_PRODUCT "Moderna COVID-19 Vaccine"
Set-DateText 2 $COL_V1_DATE    "2020-06-01"
Set-Num      2 $COL_V1_DOSE    1
Set-Text     2 $COL_V1_NOTES   "notes 1"
Set-Text     2 $COL_V2_GROUP   "COVID-19"
Set-Text     2 $COL_V2_PRODUCT "Moderna COVID-19 Vaccine"
Set-DateText 2 $COL_V2_DATE    "2020-06-20"
Set-Num      2 $COL_V2_DOSE    2
Set-Text     2 $COL_V2_NOTES   "notes 2"

# ---- Row 3 ----
Set-Text     3 $COL_V1_GROUP   "COVID-19"
Set-Text     3 $COL_V1_PRODUCT "Pfizer-BioNTech COVID-19 Vaccine"
Set-DateText 3 $COL_V1_DATE    "2020-06-02"
Set-Num      3 $COL_V1_DOSE    1
Set-Text     3 $COL_V2_GROUP   "COVID-19"
Set-Text     3 $COL_V2_PRODUCT "Pfizer-BioNTech COVID-19 Vaccine"
Set-DateText 3 $COL_V2_DATE    "2020-06-21"
Set-Num      3 $COL_V2_DOSE    2

# ---- Row 4 ----
Set-Text     4 $COL_V1_GROUP   "COVID-19"
Set-Text     4 $COL_V1_PRODUCT "Unknown"
Set-DateText 4 $COL_V1_DATE    "2020-06-04"
Set-Num      4 $COL_V1_DOSE    1
Set-Text     4 $COL_V2_GROUP   "COVID-19"
Set-Text     4 $COL_V2_PRODUCT "Unknown"
Set-DateText 4 $COL_V2_DATE    "2020-06-22"
Set-Num      4 $COL_V2_DOSE    2

# ---- Row 5 ----
Set-Text     5 $COL_V1_GROUP   "COVID-19"
Set-Text     5 $COL_V1_PRODUCT "Moderna COVID-19 Vaccine"
Set-DateText 5 $COL_V1_DATE    "2020-06-01"
Set-Num      5 $COL_V1_DOSE    1

# ---- Row 6 ----
Set-Text     6 $COL_V1_GROUP   "COVID-19"
Set-Text     6 $COL_V1_PRODUCT "Janssen (J&J) COVID-19 Vaccine"
Set-DateText 6 $COL_V1_DATE    "2020-06-03"
Set-Num      6 $COL_V1_DOSE    1

# ---- Row 7 ----
Set-Text     7 $COL_V1_GROUP   "COVID-19"
Set-Text     7 $COL_V1_PRODUCT "Unknown"
Set-DateText 7 $COL_V1_DATE    "2020-06-02"
Set-Num      7 $COL_V1_DOSE    1

# ---- Column widths (best-effort auto fit for the newly added columns) ----
$ws.Range($ws.Cells.Item(1, $COL_V1_GROUP), $ws.Cells.Item(1, $COL_V2_NOTES)).EntireColumn.AutoFit() | Out-Null

# ---- Reset the view back to the top-left corner of the sheet ----
$ws.Range("A1").Select()
